$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells (rows 46, 47, 50, 51) ---
$ws.Range("K46").Value = 8.897626167111405
$ws.Range("P46").Value = 49.29734090090081

$ws.Range("K47").Value = 8.958096993760634
$ws.Range("P47").Value = 19.47602445372739

$ws.Range("C50").Value = 7.156497794076754
$ws.Range("E50").Value = 8.791681087913137
$ws.Range("F50").Value = 9.862076629158473
$ws.Range("G50").Value = 9.082952337782311
$ws.Range("J50").Value = 4.265543567337613
$ws.Range("K50").Value = 8.357452194061691
$ws.Range("M50").Value = 6.219575926922761
$ws.Range("P50").Value = 49.60815934075665
$ws.Range("Q50").Value = 34.12762019649608

$ws.Range("F51").Value = 6.953709740148909
$ws.Range("K51").Value = 7.729722064662507
$ws.Range("P51").Value = 30.7811453972874
$ws.Range("Q51").Value = 36.95370974014891

# --- Append new rows 54-69 ---
$ws.Cells.Item(54, 1).Value = "'2025-02-14"
$ws.Cells.Item(54, 1).Style = "Normal"
$ws.Cells.Item(54, 2).Value = "abs_activity"
$arr54 = New-Object 'object[,]' 1,15
$arr54[0,0] = 10
$arr54[0,1] = 0
$arr54[0,2] = 6.734143537613756
$arr54[0,3] = 9.147439122690919
$arr54[0,4] = 9.186394141726881
$arr54[0,5] = 9.67437685676328
$arr54[0,6] = 9.451451233987699
$arr54[0,7] = 9.761738846049981
$arr54[0,8] = -11.27845954605864
$arr54[0,9] = 10
$arr54[0,10] = 8.44835319582184
$arr54[0,11] = 0
$arr54[0,12] = 0
$arr54[0,13] = 32.54188256309153
$arr54[0,14] = 38.58355482550418
$ws.Range("C54:Q54").Value = $arr54

$ws.Cells.Item(55, 1).Value = "'2025-02-14"
$ws.Cells.Item(55, 1).Style = "Normal"
$ws.Cells.Item(55, 2).Value = "rel_activity"
$arr55 = New-Object 'object[,]' 1,15
$arr55[0,0] = 10
$arr55[0,1] = 5
$arr55[0,2] = 0
$arr55[0,3] = 5.416666666666667
$arr55[0,4] = 5.601995184038527
$arr55[0,5] = 0
$arr55[0,6] = 0
$arr55[0,7] = 6.689040092207366
$arr55[0,8] = 0
$arr55[0,9] = 0
$arr55[0,10] = 6.674107142857142
$arr55[0,11] = 5
$arr55[0,12] = 5
$arr55[0,13] = 27.27610232689567
$arr55[0,14] = 22.10570675887404
$ws.Range("C55:Q55").Value = $arr55

$ws.Cells.Item(56, 1).Value = "'2025-02-14"
$ws.Cells.Item(56, 1).Style = "Normal"
$ws.Cells.Item(56, 2).Value = "abs_sleep"
$arr56 = New-Object 'object[,]' 1,15
$arr56[0,0] = 5.800000000000001
$arr56[0,1] = 9.966666666666667
$arr56[0,2] = 10
$arr56[0,3] = 6.733333333333334
$arr56[0,4] = 8.1
$arr56[0,5] = 4.333333333333332
$arr56[0,6] = 10
$arr56[0,7] = 9.333333333333334
$arr56[0,8] = 8.566666666666666
$arr56[0,9] = 10
$arr56[0,10] = 10
$arr56[0,11] = 0
$arr56[0,12] = 0
$arr56[0,13] = 52.46666666666667
$arr56[0,14] = 40.36666666666667
$ws.Range("C56:Q56").Value = $arr56

$ws.Cells.Item(57, 1).Value = "'2025-02-14"
$ws.Cells.Item(57, 1).Style = "Normal"
$ws.Cells.Item(57, 2).Value = "rel_sleep"
$arr57 = New-Object 'object[,]' 1,15
$arr57[0,0] = 0
$arr57[0,1] = 7.299645419586821
$arr57[0,2] = 10
$arr57[0,3] = 0
$arr57[0,4] = 0
$arr57[0,5] = 0
$arr57[0,6] = 7.33875338753388
$arr57[0,7] = 0
$arr57[0,8] = 0
$arr57[0,9] = 10
$arr57[0,10] = 0
$arr57[0,11] = 0
$arr57[0,12] = 0
$arr57[0,13] = 17.33875338753388
$arr57[0,14] = 17.29964541958682
$ws.Range("C57:Q57").Value = $arr57

$ws.Cells.Item(58, 1).Value = "'2025-02-15"
$ws.Cells.Item(58, 1).Style = "Normal"
$ws.Cells.Item(58, 2).Value = "abs_activity"
$arr58 = New-Object 'object[,]' 1,15
$arr58[0,0] = 6.627725557666349
$arr58[0,1] = 0
$arr58[0,2] = 9.966867400376216
$arr58[0,3] = 9.500850809362989
$arr58[0,4] = 10
$arr58[0,5] = 10
$arr58[0,6] = 10
$arr58[0,7] = 10
$arr58[0,8] = 0
$arr58[0,9] = 10
$arr58[0,10] = 9.547189050629413
$arr58[0,11] = 0
$arr58[0,12] = 0
$arr58[0,13] = 46.14178200867198
$arr58[0,14] = 39.50085080936299
$ws.Range("C58:Q58").Value = $arr58

$ws.Cells.Item(59, 1).Value = "'2025-02-15"
$ws.Cells.Item(59, 1).Style = "Normal"
$ws.Cells.Item(59, 2).Value = "rel_activity"
$arr59 = New-Object 'object[,]' 1,15
$arr59[0,0] = 0
$arr59[0,1] = 5
$arr59[0,2] = 7.109704641350211
$arr59[0,3] = 6.042957831207311
$arr59[0,4] = 10
$arr59[0,5] = 10
$arr59[0,6] = 10
$arr59[0,7] = 10
$arr59[0,8] = 0
$arr59[0,9] = 10
$arr59[0,10] = 6.494348081215556
$arr59[0,11] = 5
$arr59[0,12] = 5
$arr59[0,13] = 38.60405272256577
$arr59[0,14] = 46.04295783120731
$ws.Range("C59:Q59").Value = $arr59

$ws.Cells.Item(60, 1).Value = "'2025-02-15"
$ws.Cells.Item(60, 1).Style = "Normal"
$ws.Cells.Item(60, 2).Value = "abs_sleep"
$arr60 = New-Object 'object[,]' 1,15
$arr60[0,0] = 10
$arr60[0,1] = 10
$arr60[0,2] = 10
$arr60[0,3] = 9.699999999999999
$arr60[0,4] = 10
$arr60[0,5] = 10
$arr60[0,6] = 5.266666666666666
$arr60[0,7] = 10
$arr60[0,8] = 10
$arr60[0,9] = 9.966666666666667
$arr60[0,10] = 10
$arr60[0,11] = 0
$arr60[0,12] = 0
$arr60[0,13] = 55.26666666666667
$arr60[0,14] = 49.66666666666667
$ws.Range("C60:Q60").Value = $arr60

$ws.Cells.Item(61, 1).Value = "'2025-02-15"
$ws.Cells.Item(61, 1).Style = "Normal"
$ws.Cells.Item(61, 2).Value = "rel_sleep"
$arr61 = New-Object 'object[,]' 1,15
$arr61[0,0] = 10
$arr61[0,1] = 8.048758968553875
$arr61[0,2] = 9.083374670132347
$arr61[0,3] = 8.906411898685914
$arr61[0,4] = 9.731064618644069
$arr61[0,5] = 10
$arr61[0,6] = 0
$arr61[0,7] = 10
$arr61[0,8] = 10
$arr61[0,9] = 0
$arr61[0,10] = 0
$arr61[0,11] = 0
$arr61[0,12] = 0
$arr61[0,13] = 38.81443928877641
$arr61[0,14] = 36.95517086723979
$ws.Range("C61:Q61").Value = $arr61

$ws.Cells.Item(62, 1).Value = "'2025-02-16"
$ws.Cells.Item(62, 1).Style = "Normal"
$ws.Cells.Item(62, 2).Value = "abs_activity"
$arr62 = New-Object 'object[,]' 1,15
$arr62[0,0] = 7.306175360613179
$arr62[0,1] = 0
$arr62[0,2] = 7.686217868402406
$arr62[0,3] = 10
$arr62[0,4] = 10
$arr62[0,5] = 10
$arr62[0,6] = 7.340217078452969
$arr62[0,7] = 8.159198869261736
$arr62[0,8] = 0
$arr62[0,9] = 8.302760544226381
$arr62[0,10] = 6.521995793183863
$arr62[0,11] = 0
$arr62[0,12] = 0
$arr62[0,13] = 38.85460610065241
$arr62[0,14] = 36.46195941348812
$ws.Range("C62:Q62").Value = $arr62

$ws.Cells.Item(63, 1).Value = "'2025-02-16"
$ws.Cells.Item(63, 1).Style = "Normal"
$ws.Cells.Item(63, 2).Value = "rel_activity"
$arr63 = New-Object 'object[,]' 1,15
$arr63[0,0] = 0
$arr63[0,1] = 5
$arr63[0,2] = 0
$arr63[0,3] = 9.237934285962108
$arr63[0,4] = 9.420423127976122
$arr63[0,5] = 10
$arr63[0,6] = 0
$arr63[0,7] = 0
$arr63[0,8] = 0
$arr63[0,9] = 0
$arr63[0,10] = 0
$arr63[0,11] = 5
$arr63[0,12] = 5
$arr63[0,13] = 14.42042312797612
$arr63[0,14] = 29.23793428596211
$ws.Range("C63:Q63").Value = $arr63

$ws.Cells.Item(64, 1).Value = "'2025-02-16"
$ws.Cells.Item(64, 1).Style = "Normal"
$ws.Cells.Item(64, 2).Value = "abs_sleep"
$arr64 = New-Object 'object[,]' 1,15
$arr64[0,0] = 10
$arr64[0,1] = 9.566666666666666
$arr64[0,2] = 10
$arr64[0,3] = 2.866666666666667
$arr64[0,4] = 8.633333333333333
$arr64[0,5] = 10
$arr64[0,6] = 0
$arr64[0,7] = 10
$arr64[0,8] = 10
$arr64[0,9] = 9.800000000000001
$arr64[0,10] = 10
$arr64[0,11] = 0
$arr64[0,12] = 0
$arr64[0,13] = 48.63333333333333
$arr64[0,14] = 42.23333333333333
$ws.Range("C64:Q64").Value = $arr64

$ws.Cells.Item(65, 1).Value = "'2025-02-16"
$ws.Cells.Item(65, 1).Style = "Normal"
$ws.Cells.Item(65, 2).Value = "rel_sleep"
$arr65 = New-Object 'object[,]' 1,15
$arr65[0,0] = 9.091784989858015
$arr65[0,1] = 0
$arr65[0,2] = 0
$arr65[0,3] = 0
$arr65[0,4] = 0
$arr65[0,5] = 10
$arr65[0,6] = 0
$arr65[0,7] = 8.674692874692878
$arr65[0,8] = 8.906325948727257
$arr65[0,9] = 0
$arr65[0,10] = 0
$arr65[0,11] = 0
$arr65[0,12] = 0
$arr65[0,13] = 17.99811093858527
$arr65[0,14] = 18.67469287469288
$ws.Range("C65:Q65").Value = $arr65

$ws.Cells.Item(66, 1).Value = "'2025-02-17"
$ws.Cells.Item(66, 1).Style = "Normal"
$ws.Cells.Item(66, 2).Value = "abs_activity"
$arr66 = New-Object 'object[,]' 1,15
$arr66[0,0] = 8.165351244324929
$arr66[0,1] = 0
$arr66[0,2] = 8.677290168746785
$arr66[0,3] = 9.636502603385344
$arr66[0,4] = 8.533983577496596
$arr66[0,5] = 10
$arr66[0,6] = 9.561589654654696
$arr66[0,7] = 8.530479160535133
$arr66[0,8] = 0
$arr66[0,9] = 6.370675244861128
$arr66[0,10] = 8.413735738837513
$arr66[0,11] = 0
$arr66[0,12] = 0
$arr66[0,13] = 43.35195038406052
$arr66[0,14] = 34.53765700878161
$ws.Range("C66:Q66").Value = $arr66

$ws.Cells.Item(67, 1).Value = "'2025-02-17"
$ws.Cells.Item(67, 1).Style = "Normal"
$ws.Cells.Item(67, 2).Value = "rel_activity"
$arr67 = New-Object 'object[,]' 1,15
$arr67[0,0] = 5.644453018777788
$arr67[0,1] = 5
$arr67[0,2] = 0
$arr67[0,3] = 7.416666666666666
$arr67[0,4] = 0
$arr67[0,5] = 0
$arr67[0,6] = 6.860119047619047
$arr67[0,7] = 7.5
$arr67[0,8] = 0
$arr67[0,9] = 0
$arr67[0,10] = 6.05406746031746
$arr67[0,11] = 5
$arr67[0,12] = 5
$arr67[0,13] = 23.55863952671429
$arr67[0,14] = 24.91666666666666
$ws.Range("C67:Q67").Value = $arr67

$ws.Cells.Item(68, 1).Value = "'2025-02-17"
$ws.Cells.Item(68, 1).Style = "Normal"
$ws.Cells.Item(68, 2).Value = "abs_sleep"
$arr68 = New-Object 'object[,]' 1,15
$arr68[0,0] = 7.333333333333332
$arr68[0,1] = 9.766666666666667
$arr68[0,2] = 10
$arr68[0,3] = 10
$arr68[0,4] = 10
$arr68[0,5] = 4.800000000000001
$arr68[0,6] = 10
$arr68[0,7] = 10
$arr68[0,8] = 8.333333333333334
$arr68[0,9] = 10
$arr68[0,10] = 10
$arr68[0,11] = 0
$arr68[0,12] = 0
$arr68[0,13] = 55.66666666666666
$arr68[0,14] = 44.56666666666666
$ws.Range("C68:Q68").Value = $arr68

$ws.Cells.Item(69, 1).Value = "'2025-02-17"
$ws.Cells.Item(69, 1).Style = "Normal"
$ws.Cells.Item(69, 2).Value = "rel_sleep"
$arr69 = New-Object 'object[,]' 1,15
$arr69[0,0] = 0
$arr69[0,1] = 7
$arr69[0,2] = 7
$arr69[0,3] = 10
$arr69[0,4] = 9.234507415254239
$arr69[0,5] = 0
$arr69[0,6] = 0
$arr69[0,7] = 7.872727272727277
$arr69[0,8] = 0
$arr69[0,9] = 0
$arr69[0,10] = 7.226056815612979
$arr69[0,11] = 0
$arr69[0,12] = 0
$arr69[0,13] = 23.46056423086722
$arr69[0,14] = 24.87272727272727
$ws.Range("C69:Q69").Value = $arr69
